# Add a new "alu" result block (two lib-size rows) to the results sheet,
# matching the existing "Mult" / "Sum" / "Hamming" section layout
# (A = function name, B = lib size, P..U = synth stats + weighted total).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 is a blank spacer row; only the bold "Weight" column formatting
# carries down onto it (no value), same as the other spacer rows (7/10/11).
$ws.Range("U14").Font.Bold = $true

# Row 15: alu, lib=4
$ws.Range("A15").Value = "alu"
$ws.Range("B15").Value = 4
$ws.Range("P15").Value = 26
$ws.Range("Q15").Value = 81
$ws.Range("R15").Formula = "=SUM(P15:Q15)"
$ws.Range("S15").Value = 163
$ws.Range("T15").Formula = "=SUM(P15:R15)"
$ws.Range("U15").Formula = "=SUM(P15:Q15)+5*S15"
$ws.Range("U15").Font.Bold = $true

# Row 16: alu, lib=8
$ws.Range("B16").Value = 8
$ws.Range("P16").Value = 139
$ws.Range("Q16").Value = 220
$ws.Range("R16").Formula = "=SUM(P16:Q16)"
$ws.Range("S16").Value = 422
$ws.Range("T16").Formula = "=SUM(P16:R16)"
$ws.Range("U16").Formula = "=SUM(P16:Q16)+5*S16"
$ws.Range("U16").Font.Bold = $true

# Match the author's final cursor position recorded in the saved file.
$ws.Range("U20").Select()
